# Add a new contact row (Vehans Ayvazi) to the contacts table, matching the
# formatting pattern already used by the rest of the sheet, and hyperlink the
# new email address the same way the previous row's email is linked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12 is a fully-populated "plain text" row (Name/Phone/Email/Role all
# filled in with the default text style) - copy its formatting down into the
# new row 15 before writing values, so the new cells pick up the same style
# (s="1") used throughout the table instead of the worksheet's bare default.
$ws.Range("A12:D12").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A15").Value = "Vehans Ayvazi"
$ws.Range("B15").Value = "(818) 383-0946"
$ws.Range("C15").Value = "vehansayvazi5@gmail.com"
$ws.Range("D15").Value = "Display Data"

# Turn the new email address into a mailto hyperlink, same as was done for
# the previous row's email (C14).
$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:vehansayvazi5@gmail.com")
